$d = $word.ActiveDocument

function Rename-InlineShapeImage {
    param($range, $newName)

    $count = $range.InlineShapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $inlineShape = $range.InlineShapes($i)
        # InlineShape has no usable Name property in the Word object model -
        # it has to be converted to a (floating) Shape, renamed there, then
        # converted back to an inline picture so the drawing stays wp:inline.
        $floating = $inlineShape.ConvertToShape()
        $floating.Name = $newName
        [void]$floating.ConvertToInlineShape()
    }
}

$sec = $d.Sections(1)

# First-page header (header1.xml): BTec_Logo-Orange, image1.jpg -> image2.jpg
Rename-InlineShapeImage $sec.Headers(2).Range "image2.jpg"

# First-page footer (footer1.xml): PearsonLogo, image2.png -> image1.png
Rename-InlineShapeImage $sec.Footers(2).Range "image1.png"

# Default/primary footer (footer2.xml): PearsonLogo, image2.png -> image1.png
Rename-InlineShapeImage $sec.Footers(1).Range "image1.png"
